$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# This sheet is a daily price log for "Zapallo italiano" at the
# "Terminal Hortofruticola Agro Chillan" market. The edit inserts one
# new record by shifting the existing rows 50..131 down into rows
# 51..132 (row 131's data becomes row 132), and writes a brand-new
# record into row 49 (columns D, J, K, L, M, N, O, P, Q - the others
# - A, B, C, E, F, G, H, I, R - are constant across the whole block).
# ------------------------------------------------------------------

$lastRow = 131
$newRow = $lastRow + 1

# 1) Create the new last row by duplicating the current last row's
#    values (A..R) cell by cell, then copy the date cell's number
#    format explicitly so the new row renders the same as the rest of
#    the column (plain value-copy leaves NumberFormat at "General").
for ($c = 1; $c -le 18; $c++) {
    $ws.Cells.Item($newRow, $c).Value2 = $ws.Cells.Item($lastRow, $c).Value2
}
$ws.Cells.Item($newRow, 4).NumberFormat = $ws.Cells.Item($lastRow, 4).NumberFormat

# 2) Shift the "variable" columns down by one row, from the bottom up,
#    so every row's D/J/K/L/M/N/O/P/Q take on what used to be one row
#    above them (row 131's values already landed in row 132 via the
#    copy above, so this loop only needs to run down to row 50).
for ($r = $lastRow; $r -ge 50; $r--) {
    $src = $r - 1
    $ws.Cells.Item($r, 4).Value2  = $ws.Cells.Item($src, 4).Value2
    $ws.Cells.Item($r, 10).Value2 = $ws.Cells.Item($src, 10).Value2
    $ws.Cells.Item($r, 11).Value2 = $ws.Cells.Item($src, 11).Value2
    $ws.Cells.Item($r, 12).Value2 = $ws.Cells.Item($src, 12).Value2
    $ws.Cells.Item($r, 13).Value2 = $ws.Cells.Item($src, 13).Value2
    $ws.Cells.Item($r, 14).Value2 = $ws.Cells.Item($src, 14).Value2
    $ws.Cells.Item($r, 15).Value2 = $ws.Cells.Item($src, 15).Value2
    $ws.Cells.Item($r, 16).Value2 = $ws.Cells.Item($src, 16).Value2
    $ws.Cells.Item($r, 17).Value2 = $ws.Cells.Item($src, 17).Value2
}

# 3) Row 49 gets the brand new record values.
$ws.Cells.Item(49, 4).Value2  = 44467
$ws.Cells.Item(49, 10).Value2 = 120
$ws.Cells.Item(49, 11).Value2 = 11000
$ws.Cells.Item(49, 12).Value2 = 12000
$ws.Cells.Item(49, 13).Value2 = 11500
$ws.Cells.Item(49, 14).Value2 = "$/caja 50 unidades"
$ws.Cells.Item(49, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(49, 16).Value2 = 230
$ws.Cells.Item(49, 17).Value2 = 50
